# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker "486749 / JOSE ALEJANDRO GALTES ORDOÑEZ / 2202" row, previously
# the first data row (row 16), moves down to become the last data row
# (row 20). The worker "352956 / LUIS ALBERTO WHITEHORNE PUPO" rows move up
# to occupy rows 16-19, and their "Periodo Mora" values are reordered from
# descending (1807,1806,1805,1804) to ascending (1804,1805,1806,1807).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: now LUIS ALBERTO WHITEHORNE PUPO / 352956 / periodo 1804
$ws.Range("C16").Value = "352956"
$ws.Range("D16").Value = "LUIS ALBERTO WHITEHORNE PUPO"
$ws.Range("E16").Value = "1804"
$ws.Range("F16").Value = 31249
$ws.Range("G16").Value = 2000000

# Row 17: LUIS ALBERTO WHITEHORNE PUPO / 352956 / periodo 1805
$ws.Range("C17").Value = "352956"
$ws.Range("D17").Value = "LUIS ALBERTO WHITEHORNE PUPO"
$ws.Range("E17").Value = "1805"
$ws.Range("F17").Value = 31249
$ws.Range("G17").Value = 2000000

# Row 18: LUIS ALBERTO WHITEHORNE PUPO / 352956 / periodo 1806
$ws.Range("C18").Value = "352956"
$ws.Range("D18").Value = "LUIS ALBERTO WHITEHORNE PUPO"
$ws.Range("E18").Value = "1806"
$ws.Range("F18").Value = 31249
$ws.Range("G18").Value = 2000000

# Row 19: LUIS ALBERTO WHITEHORNE PUPO / 352956 / periodo 1807
$ws.Range("C19").Value = "352956"
$ws.Range("D19").Value = "LUIS ALBERTO WHITEHORNE PUPO"
$ws.Range("E19").Value = "1807"
$ws.Range("F19").Value = 31249
$ws.Range("G19").Value = 2000000

# Row 20: now JOSE ALEJANDRO GALTES ORDOÑEZ / 486749 / periodo 2202
$ws.Range("C20").Value = "486749"
$ws.Range("D20").Value = "JOSE ALEJANDRO GALTES ORDOÑEZ"
$ws.Range("E20").Value = "2202"
$ws.Range("F20").Value = 182000
$ws.Range("G20").Value = 4550000
